# Finished author control of theorem style classification and removal of
# italics: split a new "upright" theorem style out of TheoremStyle, and
# repoint the paragraphs that were using the plain TheoremStyle so that
# they use the new TheoremStyleUpright style instead.

$d = $word.ActiveDocument

# 1. Create the new paragraph style "TheoremStyleUpright", based on the
#    existing "TheoremStyle" custom style (wdStyleTypeParagraph = 1).
$theoremStyle = $d.Styles("TheoremStyle")
$newStyle = $d.Styles.Add("TheoremStyleUpright", 1)
$newStyle.BaseStyle = $theoremStyle
$newStyle.QuickStyle = $true

# 2. Re-point every paragraph currently formatted with "TheoremStyle" so
#    that it uses the new "TheoremStyleUpright" style instead (this is the
#    "Theorem 3.1" / "This is an inbuilt theorem." example paragraphs).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Style.NameLocal -eq "TheoremStyle") {
        $para.Style = $newStyle
    }
}

Write-Output "Applied TheoremStyle -> TheoremStyleUpright reclassification"
